$d = $word.ActiveDocument

# Locate the paragraph that contains the distinctive text
# "...3T3-L1s spreadsheet" (the last filled-in bullet under Figure 4: Lipolysis).
# The very next paragraph is the existing, empty, level-3 list item
# (ListParagraph / numId 14 / ilvl 2) after which the new bullets must be
# inserted.
$findRng = $d.Content
$found = $findRng.Find.Execute("3T3-L1s spreadsheet", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$count = $d.Paragraphs.Count
$anchorIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $findRng.Start -and $p.Range.End -ge $findRng.End) {
        $anchorIdx = $i
    }
}

$targetIdx = $anchorIdx + 1
$target = $d.Paragraphs.Item($targetIdx)

# Insert 5 new, empty list paragraphs right after the (still empty) target
# paragraph, each inheriting its ListParagraph style / numId 14 list.
$r = $target.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$r = $d.Paragraphs.Item($targetIdx + 1).Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$r = $d.Paragraphs.Item($targetIdx + 2).Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$r = $d.Paragraphs.Item($targetIdx + 3).Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$r = $d.Paragraphs.Item($targetIdx + 4).Range
$r.Collapse(0)
$r.InsertParagraphAfter()

# Paragraph 1: "In vivo lipolysis" at list level 1 (w:ilvl 0)
$p = $d.Paragraphs.Item($targetIdx + 1)
$p.Range.ListFormat.ListLevelNumber = 1
$pr = $p.Range
$pr.End = $pr.End - 1
$pr.Text = "In vivo lipolysis"

# Paragraph 2: "Rmd files" at list level 2 (w:ilvl 1)
$p = $d.Paragraphs.Item($targetIdx + 2)
$p.Range.ListFormat.ListLevelNumber = 2
$pr = $p.Range
$pr.End = $pr.End - 1
$pr.Text = "Rmd files"

# Paragraph 3: file path at list level 3 (w:ilvl 2)
$p = $d.Paragraphs.Item($targetIdx + 3)
$p.Range.ListFormat.ListLevelNumber = 3
$pr = $p.Range
$pr.End = $pr.End - 1
$pr.Text = "Harvey>Mouse work> Dexamethasone treatment> cohort A> Lipolysis folder"

# Paragraph 4: "Figures " at list level 2 (w:ilvl 1)
$p = $d.Paragraphs.Item($targetIdx + 4)
$p.Range.ListFormat.ListLevelNumber = 2
$pr = $p.Range
$pr.End = $pr.End - 1
$pr.Text = "Figures "

# Paragraph 5: file path at list level 3 (w:ilvl 2)
$p = $d.Paragraphs.Item($targetIdx + 5)
$p.Range.ListFormat.ListLevelNumber = 3
$pr = $p.Range
$pr.End = $pr.End - 1
$pr.Text = "Harvey>Mouse work> Dexamethasone treatment> cohort A> Lipolysis folder> figures folder"

Write-Output ("Inserted 5 list paragraphs after paragraph index " + $targetIdx + ". New total paragraph count: " + $d.Paragraphs.Count)
